# Add an "inhibitor value override" flag pair to the Settings sheet and
# make Settings (with the new B7 cell) the active view, matching the
# commit "Added flag to TemplateXLSX for inhibitor value override."

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# Push the existing ACCESS_CODE-and-below rows down by two so the two new
# settings land right after the EPS / INHIBITOR_MIN_VALUE / INHIBITOR_MAX_VALUE
# block.
$ws.Rows("5:6").Insert()

# Populate row 6 first so "INHIBITOR_DEFAULT_VALUE" takes the lower shared
# string slot, then row 5 with "INHIBITOR_OVERRIDE_VALUE_WITH_DEFAULT".
$ws.Range("A6").Value = "INHIBITOR_DEFAULT_VALUE"
$ws.Range("B6").Value = 0.6

$ws.Range("A5").Value = "INHIBITOR_OVERRIDE_VALUE_WITH_DEFAULT"
$ws.Range("B5").Value = $true

# The sheet is now printed in portrait orientation.
$ws.PageSetup.Orientation = 1

# Settings becomes the active/selected tab (it was Answers before), with
# B7 (the ACCESS_CODE value cell, now two rows further down) selected.
[void]$ws.Select()
[void]$ws.Range("B7").Select()
